# Refined metadata to be additional tab
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("data")

# Update the time_taken column (F) on the "data" sheet with refreshed timestamps
$ws1.Range("F2").Value = "2021-10-05 14:34:50.684962"
$ws1.Range("F3").Value = "2021-10-05 14:34:50.684970"
$ws1.Range("F4").Value = "2021-10-05 14:34:50.684973"
$ws1.Range("F5").Value = "2021-10-05 14:34:50.684976"
$ws1.Range("F6").Value = "2021-10-05 14:34:50.684979"
$ws1.Range("F7").Value = "2021-10-05 14:34:50.684982"
$ws1.Range("F8").Value = "2021-10-05 14:34:50.684985"
$ws1.Range("F9").Value = "2021-10-05 14:34:50.684987"
$ws1.Range("F10").Value = "2021-10-05 14:34:50.684990"
$ws1.Range("F11").Value = "2021-10-05 14:34:50.684993"
$ws1.Range("F12").Value = "2021-10-05 14:34:50.684995"
$ws1.Range("F13").Value = "2021-10-05 14:34:50.684998"
$ws1.Range("F14").Value = "2021-10-05 14:34:50.685000"
$ws1.Range("F15").Value = "2021-10-05 14:34:50.685003"
$ws1.Range("F16").Value = "2021-10-05 14:34:50.685006"
$ws1.Range("F17").Value = "2021-10-05 14:34:50.685008"
$ws1.Range("F18").Value = "2021-10-05 14:34:50.685011"
$ws1.Range("F19").Value = "2021-10-05 14:34:50.685014"
$ws1.Range("F20").Value = "2021-10-05 14:34:50.685016"
$ws1.Range("F21").Value = "2021-10-05 14:34:50.685019"

# Add a new "metadata" worksheet positioned after the "data" sheet
$ws2 = $wb.Worksheets.Add()
$ws2.Name = "metadata"
$ws2.Move()
# Re-fetch the worksheet reference by name: after Move() the old handle
# silently rebinds to whatever sheet now sits at the original index.
$ws2 = $wb.Worksheets.Item("metadata")

# Copy the bold/bordered/centered header style (data!B1) across the header
# row so the new cells reuse the existing style index instead of creating
# a duplicate.
$ws1.Range("B1").Copy($ws2.Range("B1:G1"))
$ws2.Range("B1").Value = "data_name"
$ws2.Range("C1").Value = "data_id"
$ws2.Range("D1").Value = "data_version"
$ws2.Range("E1").Value = "data_version_created"
$ws2.Range("F1").Value = "panel_query_time"
$ws2.Range("G1").Value = "panel_get_request"

# Copy the same styled format (data!A2) onto metadata!A2
$ws1.Range("A2").Copy($ws2.Range("A2"))
$ws2.Range("A2").Value = 0

$ws2.Range("B2").Value = "Microcephalic Primordial Dwarfism and Slender bone dysplasias"
$ws2.Range("C2").Value = 3128
$ws2.Range("D2").Value = "0.19"
$ws2.Range("E2").Value = "2021-08-18T10:10:42.996840Z"
$ws2.Range("F2").Value = "2021-10-05 14:34:50.681088"
$ws2.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/3128/?format=json"
